# Update countries & provincias Spain
# - Refresh "Datos actualizados" timestamp (07:05 -> 07:35)
# - Hungria overtakes Azerbaiyan in total cases -> rows 71/72 swap (with refreshed data)
# - Santa Lucia / Nueva Caledonia swap order (rows 197/198)
# - Seychelles / Montserrat swap order, with refreshed activos/muertes (rows 209/210)
# - Sahara Occidental / Bonaire, San Eustaquio y Saba swap order (rows 214/215)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Timestamp
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 07:35"

# Rows 71-72: Hungria moves above Azerbaiyan with updated counters
Set-Row 71 @("Hungria", 3641, 43, 1509, 1659, 0, 3, 473)
Set-Row 72 @("Azerbaiyan", 3631, 0, 2253, 1335, 0, 0, 43)

# Rows 197-198: Santa Lucia moves above Nueva Caledonia (counts unchanged)
Set-Row 197 @("Santa Lucia", 18, 0, 18, 0, 0, 0, 0)
Set-Row 198 @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)

# Rows 209-210: Montserrat moves above Seychelles with updated counters
Set-Row 209 @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
Set-Row 210 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)

# Rows 214-215: Bonaire, San Eustaquio y Saba moves above Sahara Occidental (counts unchanged)
Set-Row 214 @("Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0)
Set-Row 215 @("Sahara Occidental", 6, 0, 6, 0, 0, 0, 0)
